# "A Template for a document with tables. Compact." — retune the
# paragraph-spacing defaults baked into the Body Text / First Paragraph
# styles so body copy starts flush under a heading and First Paragraph
# gets its own breathing room above it.

$d = $word.ActiveDocument

# Body Text: keep its existing double line-spacing (line=480/auto) but
# pin space-before to 0 explicitly, instead of inheriting Normal's 12pt.
$bodyText = $d.Styles("Body Text")
$bodyText.ParagraphFormat.SpaceBefore = 0

# First Paragraph (based on Body Text): give it an explicit 12pt
# (240 twips) space-before of its own.
$firstParagraph = $d.Styles("First Paragraph")
$firstParagraph.ParagraphFormat.SpaceBefore = 12
